# rts_test.xlsx edit
# Commit: "added functions to process file data and copyText
#          getTotalRows, getTotalAmount, getInvoices, copytText button"
#
# The data-relevant part of this change (on the "Sheet1" worksheet, which is
# the active/selected tab) is the header cell G1: it used to hold the padded
# text " InvAmt " and is retitled to the clean header "InvAmt" (e.g. so it
# can be matched/used programmatically by the new getTotalAmount/getInvoices
# helpers referenced in the commit message). This introduces a brand new
# shared-string entry since the exact trimmed text didn't previously exist
# in the workbook's string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-title the InvAmt header from " InvAmt " to "InvAmt"
$ws.Range("G1").Value = "InvAmt"

# Reset the lingering stale selection (previously parked at F20) back to the
# top-left cell now that the sheet has been revisited/edited.
$ws.Range("A1").Select()
